$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Reserve type description (C16): append note about TURF
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "Type of reserve (core zone, refugio pesquero, voluntarily owned), duration of reserve (temporary or permanent), and level of protection (partial, total). If it is inside or outside a TURF"

# ---------------------------------------------------------------------------
# 2. Management plan data type (D18): Ordinal -> Binary
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = "Binary"

# ---------------------------------------------------------------------------
# 3. Row 11 ("Community knowledge of the reserves") gets struck through
#    during the meeting (superseded), and loses its yellow highlight.
# ---------------------------------------------------------------------------
$ws.Range("B11:D11").Font.Strikethrough = $true

$ws.Range("F12").Copy() | Out-Null
$ws.Range("B11:C11").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11:D11").Font.Strikethrough = $true

# ---------------------------------------------------------------------------
# 4/5. New indicators added during the meeting: Internal Regulation and
#      Perceived Effectiveness (rows 25 and 26), styled like the other
#      recently-added (highlighted) rows such as row 24.
# ---------------------------------------------------------------------------
$ws.Range("B24").Copy() | Out-Null
$ws.Range("B25:B26").PasteSpecial(-4122) | Out-Null

$ws.Range("D24").Copy() | Out-Null
$ws.Range("C25").PasteSpecial(-4122) | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(25, 2).Value = "Internal Regulation"
$ws.Cells.Item(25, 3).Value = "Does the reserve has its own regulations?"
$ws.Cells.Item(25, 4).Value = "Binary"

$ws.Cells.Item(26, 2).Value = "Perceived Effectiveness"
$ws.Cells.Item(26, 4).Value = "Binary"

# ---------------------------------------------------------------------------
# Row heights: Excel recalculated wrap-text row heights (new font metrics +
# the longer/changed text above). Match the committed values.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 26.4
$ws.Rows.Item(3).RowHeight = 26.4
$ws.Rows.Item(5).RowHeight = 26.4
$ws.Rows.Item(6).RowHeight = 26.4
$ws.Rows.Item(9).RowHeight = 145.2
$ws.Rows.Item(10).RowHeight = 26.4
$ws.Rows.Item(11).RowHeight = 79.2
$ws.Rows.Item(12).RowHeight = 26.4
$ws.Rows.Item(13).RowHeight = 26.4
$ws.Rows.Item(16).RowHeight = 39.6
$ws.Rows.Item(18).RowHeight = 39.6
$ws.Rows.Item(20).RowHeight = 26.4
$ws.Rows.Item(22).RowHeight = 39.6
$ws.Rows.Item(23).RowHeight = 39.6
$ws.Rows.Item(24).RowHeight = 26.4

$ws.Range("C26").Select()
